# -----------------------------------------------------------------------
# Commit: "feat: add 2022-Q4 data"
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (shifting the
#    other quarterly sheets one position to the right).
# 2. Populate "总计" (summary) row 2 with the new 2022-Q4 totals and
#    push the previously-existing rows down by one, adding a brand new
#    trailing row for 2021-Q1.
# 3. Fill the new "2022-Q4" sheet with its per-fund breakdown rows.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Insert + rename the new sheet -------------------------------
$afterTotal = $wb.Worksheets.Item(2)          # current "2022-Q3" sheet
$wsQ4 = $wb.Worksheets.Add($afterTotal)       # new sheet inserted before it
$wsQ4.Name = "2022-Q4"

# ---- 2. Update the "总计" summary sheet ------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Cells.Item(2,2).Value = "2022-Q4"
$wsTotal.Cells.Item(2,3).Value = 29
$wsTotal.Cells.Item(2,4).Value = 11.6

$wsTotal.Cells.Item(3,2).Value = "2022-Q3"
$wsTotal.Cells.Item(3,3).Value = 30
$wsTotal.Cells.Item(3,4).Value = 12.67

$wsTotal.Cells.Item(4,2).Value = "2022-Q2"
$wsTotal.Cells.Item(4,3).Value = 23
$wsTotal.Cells.Item(4,4).Value = 12.05

$wsTotal.Cells.Item(5,2).Value = "2022-Q1"
$wsTotal.Cells.Item(5,3).Value = 27
$wsTotal.Cells.Item(5,4).Value = 24.25

$wsTotal.Cells.Item(6,2).Value = "2021-Q4"
$wsTotal.Cells.Item(6,3).Value = 23
$wsTotal.Cells.Item(6,4).Value = 23.35

$wsTotal.Cells.Item(7,2).Value = "2021-Q2"
$wsTotal.Cells.Item(7,3).Value = 24
$wsTotal.Cells.Item(7,4).Value = 18.06

$wsTotal.Cells.Item(8,1).Value = 6
$wsTotal.Cells.Item(8,2).Value = "2021-Q1"
$wsTotal.Cells.Item(8,3).Value = 23
$wsTotal.Cells.Item(8,4).Value = 15.26

# Column A on row 8 is new -- copy the number-style (bold/bordered) of
# the existing A column cells onto it so it matches its neighbours.
$wsTotal.Cells.Item(7,1).Copy()
$wsTotal.Cells.Item(8,1).PasteSpecial(-4122)

# ---- 3. Populate the new "2022-Q4" sheet -----------------------------
$wsQ4.Cells.Item(1,2).Value = "'基金代码"
$wsQ4.Cells.Item(1,2).Style = "Normal"
$wsQ4.Cells.Item(1,3).Value = "'基金名称"
$wsQ4.Cells.Item(1,3).Style = "Normal"
$wsQ4.Cells.Item(1,4).Value = "'基金规模"
$wsQ4.Cells.Item(1,4).Style = "Normal"
$wsQ4.Cells.Item(1,5).Value = "'股票总仓位"
$wsQ4.Cells.Item(1,5).Style = "Normal"
$wsQ4.Cells.Item(1,6).Value = "'仓位占比"
$wsQ4.Cells.Item(1,6).Style = "Normal"
$wsQ4.Cells.Item(1,7).Value = "'持有市值(亿元)"
$wsQ4.Cells.Item(1,7).Style = "Normal"
$wsQ4.Cells.Item(1,8).Value = "'仓位排名"
$wsQ4.Cells.Item(1,8).Style = "Normal"
$wsQ4.Cells.Item(2,1).Value = 0
$wsQ4.Cells.Item(2,2).Value = "'159941"
$wsQ4.Cells.Item(2,2).Style = "Normal"
$wsQ4.Cells.Item(2,3).Value = "'广发纳斯达克100ETF（QDII）"
$wsQ4.Cells.Item(2,3).Style = "Normal"
$wsQ4.Cells.Item(2,4).Value = "'114.77"
$wsQ4.Cells.Item(2,4).Style = "Normal"
$wsQ4.Cells.Item(2,5).Value = "'90.42"
$wsQ4.Cells.Item(2,5).Style = "Normal"
$wsQ4.Cells.Item(2,6).Value = "'3.37"
$wsQ4.Cells.Item(2,6).Style = "Normal"
$wsQ4.Cells.Item(2,7).Value = "'3.8677"
$wsQ4.Cells.Item(2,7).Style = "Normal"
$wsQ4.Cells.Item(2,8).Value = 5
$wsQ4.Cells.Item(3,1).Value = 1
$wsQ4.Cells.Item(3,2).Value = "'513100"
$wsQ4.Cells.Item(3,2).Style = "Normal"
$wsQ4.Cells.Item(3,3).Value = "'国泰纳斯达克100（QDII-ETF）"
$wsQ4.Cells.Item(3,3).Style = "Normal"
$wsQ4.Cells.Item(3,4).Value = "'51.50"
$wsQ4.Cells.Item(3,4).Style = "Normal"
$wsQ4.Cells.Item(3,5).Value = "'90.72"
$wsQ4.Cells.Item(3,5).Style = "Normal"
$wsQ4.Cells.Item(3,6).Value = "'2.89"
$wsQ4.Cells.Item(3,6).Style = "Normal"
$wsQ4.Cells.Item(3,7).Value = "'1.4884"
$wsQ4.Cells.Item(3,7).Style = "Normal"
$wsQ4.Cells.Item(3,8).Value = 5
$wsQ4.Cells.Item(4,1).Value = 2
$wsQ4.Cells.Item(4,2).Value = "'513500"
$wsQ4.Cells.Item(4,2).Style = "Normal"
$wsQ4.Cells.Item(4,3).Value = "'博时标普500ETF（QDII）"
$wsQ4.Cells.Item(4,3).Style = "Normal"
$wsQ4.Cells.Item(4,4).Value = "'93.59"
$wsQ4.Cells.Item(4,4).Style = "Normal"
$wsQ4.Cells.Item(4,5).Value = "'95.36"
$wsQ4.Cells.Item(4,5).Style = "Normal"
$wsQ4.Cells.Item(4,6).Value = "'1.56"
$wsQ4.Cells.Item(4,6).Style = "Normal"
$wsQ4.Cells.Item(4,7).Value = "'1.4600"
$wsQ4.Cells.Item(4,7).Style = "Normal"
$wsQ4.Cells.Item(4,8).Value = 5
$wsQ4.Cells.Item(5,1).Value = 3
$wsQ4.Cells.Item(5,2).Value = "'159632"
$wsQ4.Cells.Item(5,2).Style = "Normal"
$wsQ4.Cells.Item(5,3).Value = "'华安纳斯达克100ETF（QDII）"
$wsQ4.Cells.Item(5,3).Style = "Normal"
$wsQ4.Cells.Item(5,4).Value = "'30.39"
$wsQ4.Cells.Item(5,4).Style = "Normal"
$wsQ4.Cells.Item(5,5).Value = "'87.77"
$wsQ4.Cells.Item(5,5).Style = "Normal"
$wsQ4.Cells.Item(5,6).Value = "'3.34"
$wsQ4.Cells.Item(5,6).Style = "Normal"
$wsQ4.Cells.Item(5,7).Value = "'1.0150"
$wsQ4.Cells.Item(5,7).Style = "Normal"
$wsQ4.Cells.Item(5,8).Value = 5
$wsQ4.Cells.Item(6,1).Value = 4
$wsQ4.Cells.Item(6,2).Value = "'000834"
$wsQ4.Cells.Item(6,2).Style = "Normal"
$wsQ4.Cells.Item(6,3).Value = "'大成纳斯达克100指数（QDII）"
$wsQ4.Cells.Item(6,3).Style = "Normal"
$wsQ4.Cells.Item(6,4).Value = "'15.61"
$wsQ4.Cells.Item(6,4).Style = "Normal"
$wsQ4.Cells.Item(6,5).Value = "'81.77"
$wsQ4.Cells.Item(6,5).Style = "Normal"
$wsQ4.Cells.Item(6,6).Value = "'3.64"
$wsQ4.Cells.Item(6,6).Style = "Normal"
$wsQ4.Cells.Item(6,7).Value = "'0.5682"
$wsQ4.Cells.Item(6,7).Style = "Normal"
$wsQ4.Cells.Item(6,8).Value = 4
$wsQ4.Cells.Item(7,1).Value = 5
$wsQ4.Cells.Item(7,2).Value = "'160213"
$wsQ4.Cells.Item(7,2).Style = "Normal"
$wsQ4.Cells.Item(7,3).Value = "'国泰纳斯达克100指数（QDII）"
$wsQ4.Cells.Item(7,3).Style = "Normal"
$wsQ4.Cells.Item(7,4).Value = "'15.65"
$wsQ4.Cells.Item(7,4).Style = "Normal"
$wsQ4.Cells.Item(7,5).Value = "'90.80"
$wsQ4.Cells.Item(7,5).Style = "Normal"
$wsQ4.Cells.Item(7,6).Value = "'3.45"
$wsQ4.Cells.Item(7,6).Style = "Normal"
$wsQ4.Cells.Item(7,7).Value = "'0.5399"
$wsQ4.Cells.Item(7,7).Style = "Normal"
$wsQ4.Cells.Item(7,8).Value = 5
$wsQ4.Cells.Item(8,1).Value = 6
$wsQ4.Cells.Item(8,2).Value = "'000043"
$wsQ4.Cells.Item(8,2).Style = "Normal"
$wsQ4.Cells.Item(8,3).Value = "'嘉实美国成长股票（QDII）人民币"
$wsQ4.Cells.Item(8,3).Style = "Normal"
$wsQ4.Cells.Item(8,4).Value = "'12.69"
$wsQ4.Cells.Item(8,4).Style = "Normal"
$wsQ4.Cells.Item(8,5).Value = "'92.23"
$wsQ4.Cells.Item(8,5).Style = "Normal"
$wsQ4.Cells.Item(8,6).Value = "'3.90"
$wsQ4.Cells.Item(8,6).Style = "Normal"
$wsQ4.Cells.Item(8,7).Value = "'0.4949"
$wsQ4.Cells.Item(8,7).Style = "Normal"
$wsQ4.Cells.Item(8,8).Value = 3
$wsQ4.Cells.Item(9,1).Value = 7
$wsQ4.Cells.Item(9,2).Value = "'000044"
$wsQ4.Cells.Item(9,2).Style = "Normal"
$wsQ4.Cells.Item(9,3).Value = "'嘉实美国成长股票（QDII）美元现汇"
$wsQ4.Cells.Item(9,3).Style = "Normal"
$wsQ4.Cells.Item(9,4).Value = "'12.69"
$wsQ4.Cells.Item(9,4).Style = "Normal"
$wsQ4.Cells.Item(9,5).Value = "'92.23"
$wsQ4.Cells.Item(9,5).Style = "Normal"
$wsQ4.Cells.Item(9,6).Value = "'3.90"
$wsQ4.Cells.Item(9,6).Style = "Normal"
$wsQ4.Cells.Item(9,7).Value = "'0.4949"
$wsQ4.Cells.Item(9,7).Style = "Normal"
$wsQ4.Cells.Item(9,8).Value = 3
$wsQ4.Cells.Item(10,1).Value = 8
$wsQ4.Cells.Item(10,2).Value = "'513300"
$wsQ4.Cells.Item(10,2).Style = "Normal"
$wsQ4.Cells.Item(10,3).Value = "'华夏纳斯达克100ETF（QDII）"
$wsQ4.Cells.Item(10,3).Style = "Normal"
$wsQ4.Cells.Item(10,4).Value = "'12.43"
$wsQ4.Cells.Item(10,4).Style = "Normal"
$wsQ4.Cells.Item(10,5).Value = "'97.54"
$wsQ4.Cells.Item(10,5).Style = "Normal"
$wsQ4.Cells.Item(10,6).Value = "'3.71"
$wsQ4.Cells.Item(10,6).Style = "Normal"
$wsQ4.Cells.Item(10,7).Value = "'0.4612"
$wsQ4.Cells.Item(10,7).Style = "Normal"
$wsQ4.Cells.Item(10,8).Value = 4
$wsQ4.Cells.Item(11,1).Value = 9
$wsQ4.Cells.Item(11,2).Value = "'161130"
$wsQ4.Cells.Item(11,2).Style = "Normal"
$wsQ4.Cells.Item(11,3).Value = "'易方达纳斯达克100指数人民币（QDII-LOF）"
$wsQ4.Cells.Item(11,3).Style = "Normal"
$wsQ4.Cells.Item(11,4).Value = "'7.77"
$wsQ4.Cells.Item(11,4).Style = "Normal"
$wsQ4.Cells.Item(11,5).Value = "'90.34"
$wsQ4.Cells.Item(11,5).Style = "Normal"
$wsQ4.Cells.Item(11,6).Value = "'3.43"
$wsQ4.Cells.Item(11,6).Style = "Normal"
$wsQ4.Cells.Item(11,7).Value = "'0.2665"
$wsQ4.Cells.Item(11,7).Style = "Normal"
$wsQ4.Cells.Item(11,8).Value = 5
$wsQ4.Cells.Item(12,1).Value = 10
$wsQ4.Cells.Item(12,2).Value = "'003722"
$wsQ4.Cells.Item(12,2).Style = "Normal"
$wsQ4.Cells.Item(12,3).Value = "'易方达纳斯达克100指数美元（QDII-LOF）A"
$wsQ4.Cells.Item(12,3).Style = "Normal"
$wsQ4.Cells.Item(12,4).Value = "'7.77"
$wsQ4.Cells.Item(12,4).Style = "Normal"
$wsQ4.Cells.Item(12,5).Value = "'90.34"
$wsQ4.Cells.Item(12,5).Style = "Normal"
$wsQ4.Cells.Item(12,6).Value = "'3.43"
$wsQ4.Cells.Item(12,6).Style = "Normal"
$wsQ4.Cells.Item(12,7).Value = "'0.2665"
$wsQ4.Cells.Item(12,7).Style = "Normal"
$wsQ4.Cells.Item(12,8).Value = 5
$wsQ4.Cells.Item(13,1).Value = 11
$wsQ4.Cells.Item(13,2).Value = "'012860"
$wsQ4.Cells.Item(13,2).Style = "Normal"
$wsQ4.Cells.Item(13,3).Value = "'易方达标普500指数（QDII-LOF）人民币 C"
$wsQ4.Cells.Item(13,3).Style = "Normal"
$wsQ4.Cells.Item(13,4).Value = "'4.75"
$wsQ4.Cells.Item(13,4).Style = "Normal"
$wsQ4.Cells.Item(13,5).Value = "'91.65"
$wsQ4.Cells.Item(13,5).Style = "Normal"
$wsQ4.Cells.Item(13,6).Value = "'1.50"
$wsQ4.Cells.Item(13,6).Style = "Normal"
$wsQ4.Cells.Item(13,7).Value = "'0.0712"
$wsQ4.Cells.Item(13,7).Style = "Normal"
$wsQ4.Cells.Item(13,8).Value = 5
$wsQ4.Cells.Item(14,1).Value = 12
$wsQ4.Cells.Item(14,2).Value = "'161125"
$wsQ4.Cells.Item(14,2).Style = "Normal"
$wsQ4.Cells.Item(14,3).Value = "'易方达标普500指数（QDII-LOF）人民币"
$wsQ4.Cells.Item(14,3).Style = "Normal"
$wsQ4.Cells.Item(14,4).Value = "'4.75"
$wsQ4.Cells.Item(14,4).Style = "Normal"
$wsQ4.Cells.Item(14,5).Value = "'91.65"
$wsQ4.Cells.Item(14,5).Style = "Normal"
$wsQ4.Cells.Item(14,6).Value = "'1.50"
$wsQ4.Cells.Item(14,6).Style = "Normal"
$wsQ4.Cells.Item(14,7).Value = "'0.0712"
$wsQ4.Cells.Item(14,7).Style = "Normal"
$wsQ4.Cells.Item(14,8).Value = 5
$wsQ4.Cells.Item(15,1).Value = 13
$wsQ4.Cells.Item(15,2).Value = "'003718"
$wsQ4.Cells.Item(15,2).Style = "Normal"
$wsQ4.Cells.Item(15,3).Value = "'易方达标普500指数（QDII-LOF）美元A"
$wsQ4.Cells.Item(15,3).Style = "Normal"
$wsQ4.Cells.Item(15,4).Value = "'4.65"
$wsQ4.Cells.Item(15,4).Style = "Normal"
$wsQ4.Cells.Item(15,5).Value = "'91.65"
$wsQ4.Cells.Item(15,5).Style = "Normal"
$wsQ4.Cells.Item(15,6).Value = "'1.50"
$wsQ4.Cells.Item(15,6).Style = "Normal"
$wsQ4.Cells.Item(15,7).Value = "'0.0698"
$wsQ4.Cells.Item(15,7).Style = "Normal"
$wsQ4.Cells.Item(15,8).Value = 5
$wsQ4.Cells.Item(16,1).Value = 14
$wsQ4.Cells.Item(16,2).Value = "'006792"
$wsQ4.Cells.Item(16,2).Style = "Normal"
$wsQ4.Cells.Item(16,3).Value = "'鹏华香港美国互联网股票（LOF）美元现汇"
$wsQ4.Cells.Item(16,3).Style = "Normal"
$wsQ4.Cells.Item(16,4).Value = "'1.29"
$wsQ4.Cells.Item(16,4).Style = "Normal"
$wsQ4.Cells.Item(16,5).Value = "'88.46"
$wsQ4.Cells.Item(16,5).Style = "Normal"
$wsQ4.Cells.Item(16,6).Value = "'4.77"
$wsQ4.Cells.Item(16,6).Style = "Normal"
$wsQ4.Cells.Item(16,7).Value = "'0.0615"
$wsQ4.Cells.Item(16,7).Style = "Normal"
$wsQ4.Cells.Item(16,8).Value = 6
$wsQ4.Cells.Item(17,1).Value = 15
$wsQ4.Cells.Item(17,2).Value = "'160644"
$wsQ4.Cells.Item(17,2).Style = "Normal"
$wsQ4.Cells.Item(17,3).Value = "'鹏华香港美国互联网股票（LOF）人民币"
$wsQ4.Cells.Item(17,3).Style = "Normal"
$wsQ4.Cells.Item(17,4).Value = "'1.29"
$wsQ4.Cells.Item(17,4).Style = "Normal"
$wsQ4.Cells.Item(17,5).Value = "'88.46"
$wsQ4.Cells.Item(17,5).Style = "Normal"
$wsQ4.Cells.Item(17,6).Value = "'4.77"
$wsQ4.Cells.Item(17,6).Style = "Normal"
$wsQ4.Cells.Item(17,7).Value = "'0.0615"
$wsQ4.Cells.Item(17,7).Style = "Normal"
$wsQ4.Cells.Item(17,8).Value = 6
$wsQ4.Cells.Item(18,1).Value = 16
$wsQ4.Cells.Item(18,2).Value = "'016532"
$wsQ4.Cells.Item(18,2).Style = "Normal"
$wsQ4.Cells.Item(18,3).Value = "'嘉实纳斯达克100指数（QDII）A人民币"
$wsQ4.Cells.Item(18,3).Style = "Normal"
$wsQ4.Cells.Item(18,4).Value = "'1.12"
$wsQ4.Cells.Item(18,4).Style = "Normal"
$wsQ4.Cells.Item(18,5).Value = "'94.67"
$wsQ4.Cells.Item(18,5).Style = "Normal"
$wsQ4.Cells.Item(18,6).Value = "'3.60"
$wsQ4.Cells.Item(18,6).Style = "Normal"
$wsQ4.Cells.Item(18,7).Value = "'0.0403"
$wsQ4.Cells.Item(18,7).Style = "Normal"
$wsQ4.Cells.Item(18,8).Value = 5
$wsQ4.Cells.Item(19,1).Value = 17
$wsQ4.Cells.Item(19,2).Value = "'016533"
$wsQ4.Cells.Item(19,2).Style = "Normal"
$wsQ4.Cells.Item(19,3).Value = "'嘉实纳斯达克100指数（QDII）C人民币"
$wsQ4.Cells.Item(19,3).Style = "Normal"
$wsQ4.Cells.Item(19,4).Value = "'1.12"
$wsQ4.Cells.Item(19,4).Style = "Normal"
$wsQ4.Cells.Item(19,5).Value = "'94.67"
$wsQ4.Cells.Item(19,5).Style = "Normal"
$wsQ4.Cells.Item(19,6).Value = "'3.60"
$wsQ4.Cells.Item(19,6).Style = "Normal"
$wsQ4.Cells.Item(19,7).Value = "'0.0403"
$wsQ4.Cells.Item(19,7).Style = "Normal"
$wsQ4.Cells.Item(19,8).Value = 5
$wsQ4.Cells.Item(20,1).Value = 18
$wsQ4.Cells.Item(20,2).Value = "'016534"
$wsQ4.Cells.Item(20,2).Style = "Normal"
$wsQ4.Cells.Item(20,3).Value = "'嘉实纳斯达克100指数（QDII）A美元现汇"
$wsQ4.Cells.Item(20,3).Style = "Normal"
$wsQ4.Cells.Item(20,4).Value = "'1.12"
$wsQ4.Cells.Item(20,4).Style = "Normal"
$wsQ4.Cells.Item(20,5).Value = "'94.67"
$wsQ4.Cells.Item(20,5).Style = "Normal"
$wsQ4.Cells.Item(20,6).Value = "'3.60"
$wsQ4.Cells.Item(20,6).Style = "Normal"
$wsQ4.Cells.Item(20,7).Value = "'0.0403"
$wsQ4.Cells.Item(20,7).Style = "Normal"
$wsQ4.Cells.Item(20,8).Value = 5
$wsQ4.Cells.Item(21,1).Value = 19
$wsQ4.Cells.Item(21,2).Value = "'016535"
$wsQ4.Cells.Item(21,2).Style = "Normal"
$wsQ4.Cells.Item(21,3).Value = "'嘉实纳斯达克100指数（QDII）C美元现汇"
$wsQ4.Cells.Item(21,3).Style = "Normal"
$wsQ4.Cells.Item(21,4).Value = "'1.12"
$wsQ4.Cells.Item(21,4).Style = "Normal"
$wsQ4.Cells.Item(21,5).Value = "'94.67"
$wsQ4.Cells.Item(21,5).Style = "Normal"
$wsQ4.Cells.Item(21,6).Value = "'3.60"
$wsQ4.Cells.Item(21,6).Style = "Normal"
$wsQ4.Cells.Item(21,7).Value = "'0.0403"
$wsQ4.Cells.Item(21,7).Style = "Normal"
$wsQ4.Cells.Item(21,8).Value = 5
$wsQ4.Cells.Item(22,1).Value = 20
$wsQ4.Cells.Item(22,2).Value = "'016055"
$wsQ4.Cells.Item(22,2).Style = "Normal"
$wsQ4.Cells.Item(22,3).Value = "'博时纳斯达克100指数（QDII）A人民币"
$wsQ4.Cells.Item(22,3).Style = "Normal"
$wsQ4.Cells.Item(22,4).Value = "'1.06"
$wsQ4.Cells.Item(22,4).Style = "Normal"
$wsQ4.Cells.Item(22,5).Value = "'90.62"
$wsQ4.Cells.Item(22,5).Style = "Normal"
$wsQ4.Cells.Item(22,6).Value = "'3.44"
$wsQ4.Cells.Item(22,6).Style = "Normal"
$wsQ4.Cells.Item(22,7).Value = "'0.0365"
$wsQ4.Cells.Item(22,7).Style = "Normal"
$wsQ4.Cells.Item(22,8).Value = 5
$wsQ4.Cells.Item(23,1).Value = 21
$wsQ4.Cells.Item(23,2).Value = "'016057"
$wsQ4.Cells.Item(23,2).Style = "Normal"
$wsQ4.Cells.Item(23,3).Value = "'博时纳斯达克100指数（QDII）C人民币"
$wsQ4.Cells.Item(23,3).Style = "Normal"
$wsQ4.Cells.Item(23,4).Value = "'1.06"
$wsQ4.Cells.Item(23,4).Style = "Normal"
$wsQ4.Cells.Item(23,5).Value = "'90.62"
$wsQ4.Cells.Item(23,5).Style = "Normal"
$wsQ4.Cells.Item(23,6).Value = "'3.44"
$wsQ4.Cells.Item(23,6).Style = "Normal"
$wsQ4.Cells.Item(23,7).Value = "'0.0365"
$wsQ4.Cells.Item(23,7).Style = "Normal"
$wsQ4.Cells.Item(23,8).Value = 5
$wsQ4.Cells.Item(24,1).Value = 22
$wsQ4.Cells.Item(24,2).Value = "'016056"
$wsQ4.Cells.Item(24,2).Style = "Normal"
$wsQ4.Cells.Item(24,3).Value = "'博时纳斯达克100指数（QDII）A美元现汇"
$wsQ4.Cells.Item(24,3).Style = "Normal"
$wsQ4.Cells.Item(24,4).Value = "'1.06"
$wsQ4.Cells.Item(24,4).Style = "Normal"
$wsQ4.Cells.Item(24,5).Value = "'90.62"
$wsQ4.Cells.Item(24,5).Style = "Normal"
$wsQ4.Cells.Item(24,6).Value = "'3.44"
$wsQ4.Cells.Item(24,6).Style = "Normal"
$wsQ4.Cells.Item(24,7).Value = "'0.0365"
$wsQ4.Cells.Item(24,7).Style = "Normal"
$wsQ4.Cells.Item(24,8).Value = 5
$wsQ4.Cells.Item(25,1).Value = 23
$wsQ4.Cells.Item(25,2).Value = "'016058"
$wsQ4.Cells.Item(25,2).Style = "Normal"
$wsQ4.Cells.Item(25,3).Value = "'博时纳斯达克100指数（QDII）C美元现汇"
$wsQ4.Cells.Item(25,3).Style = "Normal"
$wsQ4.Cells.Item(25,4).Value = "'1.06"
$wsQ4.Cells.Item(25,4).Style = "Normal"
$wsQ4.Cells.Item(25,5).Value = "'90.62"
$wsQ4.Cells.Item(25,5).Style = "Normal"
$wsQ4.Cells.Item(25,6).Value = "'3.44"
$wsQ4.Cells.Item(25,6).Style = "Normal"
$wsQ4.Cells.Item(25,7).Value = "'0.0365"
$wsQ4.Cells.Item(25,7).Style = "Normal"
$wsQ4.Cells.Item(25,8).Value = 5
$wsQ4.Cells.Item(26,1).Value = 24
$wsQ4.Cells.Item(26,2).Value = "'159612"
$wsQ4.Cells.Item(26,2).Style = "Normal"
$wsQ4.Cells.Item(26,3).Value = "'国泰标普500ETF（QDII）"
$wsQ4.Cells.Item(26,3).Style = "Normal"
$wsQ4.Cells.Item(26,4).Value = "'0.86"
$wsQ4.Cells.Item(26,4).Style = "Normal"
$wsQ4.Cells.Item(26,5).Value = "'94.21"
$wsQ4.Cells.Item(26,5).Style = "Normal"
$wsQ4.Cells.Item(26,6).Value = "'1.55"
$wsQ4.Cells.Item(26,6).Style = "Normal"
$wsQ4.Cells.Item(26,7).Value = "'0.0133"
$wsQ4.Cells.Item(26,7).Style = "Normal"
$wsQ4.Cells.Item(26,8).Value = 5
$wsQ4.Cells.Item(27,1).Value = 25
$wsQ4.Cells.Item(27,2).Value = "'012870"
$wsQ4.Cells.Item(27,2).Style = "Normal"
$wsQ4.Cells.Item(27,3).Value = "'易方达纳斯达克100指数人民币（QDII-LOF）C"
$wsQ4.Cells.Item(27,3).Style = "Normal"
$wsQ4.Cells.Item(27,4).Value = "'0.21"
$wsQ4.Cells.Item(27,4).Style = "Normal"
$wsQ4.Cells.Item(27,5).Value = "'90.34"
$wsQ4.Cells.Item(27,5).Style = "Normal"
$wsQ4.Cells.Item(27,6).Value = "'3.43"
$wsQ4.Cells.Item(27,6).Style = "Normal"
$wsQ4.Cells.Item(27,7).Value = "'0.0072"
$wsQ4.Cells.Item(27,7).Style = "Normal"
$wsQ4.Cells.Item(27,8).Value = 5
$wsQ4.Cells.Item(28,1).Value = 26
$wsQ4.Cells.Item(28,2).Value = "'012871"
$wsQ4.Cells.Item(28,2).Style = "Normal"
$wsQ4.Cells.Item(28,3).Value = "'易方达纳斯达克100指数美元（QDII-LOF）C"
$wsQ4.Cells.Item(28,3).Style = "Normal"
$wsQ4.Cells.Item(28,4).Value = "'0.21"
$wsQ4.Cells.Item(28,4).Style = "Normal"
$wsQ4.Cells.Item(28,5).Value = "'90.34"
$wsQ4.Cells.Item(28,5).Style = "Normal"
$wsQ4.Cells.Item(28,6).Value = "'3.43"
$wsQ4.Cells.Item(28,6).Style = "Normal"
$wsQ4.Cells.Item(28,7).Value = "'0.0072"
$wsQ4.Cells.Item(28,7).Style = "Normal"
$wsQ4.Cells.Item(28,8).Value = 5
$wsQ4.Cells.Item(29,1).Value = 27
$wsQ4.Cells.Item(29,2).Value = "'159655"
$wsQ4.Cells.Item(29,2).Style = "Normal"
$wsQ4.Cells.Item(29,3).Value = "'华夏标普500ETF（QDII）"
$wsQ4.Cells.Item(29,3).Style = "Normal"
$wsQ4.Cells.Item(29,4).Value = "'0.21"
$wsQ4.Cells.Item(29,4).Style = "Normal"
$wsQ4.Cells.Item(29,5).Value = "'93.70"
$wsQ4.Cells.Item(29,5).Style = "Normal"
$wsQ4.Cells.Item(29,6).Value = "'1.54"
$wsQ4.Cells.Item(29,6).Style = "Normal"
$wsQ4.Cells.Item(29,7).Value = "'0.0032"
$wsQ4.Cells.Item(29,7).Style = "Normal"
$wsQ4.Cells.Item(29,8).Value = 4
$wsQ4.Cells.Item(30,1).Value = 28
$wsQ4.Cells.Item(30,2).Value = "'012861"
$wsQ4.Cells.Item(30,2).Style = "Normal"
$wsQ4.Cells.Item(30,3).Value = "'易方达标普500指数（QDII-LOF）美元 C"
$wsQ4.Cells.Item(30,3).Style = "Normal"
$wsQ4.Cells.Item(30,4).Value = "'0.10"
$wsQ4.Cells.Item(30,4).Style = "Normal"
$wsQ4.Cells.Item(30,5).Value = "'91.65"
$wsQ4.Cells.Item(30,5).Style = "Normal"
$wsQ4.Cells.Item(30,6).Value = "'1.50"
$wsQ4.Cells.Item(30,6).Style = "Normal"
$wsQ4.Cells.Item(30,7).Value = "'0.0015"
$wsQ4.Cells.Item(30,7).Style = "Normal"
$wsQ4.Cells.Item(30,8).Value = 5

# Re-apply the workbook's header style (bold + bordered, centred) to the
# "2022-Q4" header row and its row-index column, matching the other
# quarterly sheets (copied from the "总计" sheet, which already carries
# that exact style).
$wsTotal.Cells.Item(2,1).Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsQ4.Range("A2:A30").PasteSpecial(-4122)
